# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For every data row, the comma-separated list of recorders in column G
# is reversed in order (e.g. "a@x.com, System" -> "System, a@x.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ', '
        if ($parts.Count -gt 1) {
            $rev = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $rev += $parts[$i]
            }
            $cell.Value = $rev -join ', '
        }
    }
}
